# Insert a new row at position 24. Excel shifts rows 24:80 down to 25:81,
# which reproduces the "old row N becomes new row N+1" pattern seen in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(24).Insert()

# The newly inserted row 24 is blank; populate it with the same data that was
# in the original row 24 (now shifted to row 25), then correct only the date.
$ws.Range("A24:T24").Value2 = $ws.Range("A25:T25").Value2
$ws.Range("D24").Value2 = 44519
